# The deck currently carries two themes:
#   theme1.xml = "Office Theme" (clrScheme "Office")      -- unused/orphaned by the active master
#   theme2.xml = "Integral"     (clrScheme "Red Violet")  -- the theme actually applied to
#                                                             the slide master / presentation
#
# The target edit swaps their contents so the presentation's active theme
# (theme2.xml, reached via $p.SlideMaster) becomes the plain "Office Theme"
# color scheme instead of "Integral"/"Red Violet". We reproduce that by
# rewriting each of the twelve theme colors on the slide master's
# ColorScheme, in clrScheme order: dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink -- using the exact RGB values the "Office Theme" used.
#
# (PowerPoint's RGB() helper is not available in this host, so the BGR-packed
# long values -- R + G*256 + B*65536, exactly what RGB() would return -- are
# supplied directly.)

$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$cs = $m.ColorScheme

$cs.Colors(1).RGB  = 0        # dk1     000000
$cs.Colors(2).RGB  = 16777215 # lt1     FFFFFF
$cs.Colors(3).RGB  = 6968388  # dk2     44546A
$cs.Colors(4).RGB  = 15132391 # lt2     E7E6E6
$cs.Colors(5).RGB  = 13998939 # accent1 5B9BD5
$cs.Colors(6).RGB  = 3243501  # accent2 ED7D31
$cs.Colors(7).RGB  = 10855845 # accent3 A5A5A5
$cs.Colors(8).RGB  = 49407    # accent4 FFC000
$cs.Colors(9).RGB  = 12874308 # accent5 4472C4
$cs.Colors(10).RGB = 4697456  # accent6 70AD47
$cs.Colors(11).RGB = 12673797 # hlink   0563C1
$cs.Colors(12).RGB = 7491477  # folHlink 954F72
